# Editar Dados Cadastrais para correção
# Replace the old "Tela ESCOLARIDADE" / "Tela Proficiência" block (rows 79-80)
# with the new "Idioma" / "Órgão Financiador" / "UF Nascimento" review notes
# that now span rows 79-87 on the "Cadastrar Pesquisa" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cadastrar Pesquisa")

$ws.Range("B79").Value = "Idioma"
$ws.Range("C79").Value = "Corrigir insere linhas (iniciar na linha 1)"

$ws.Range("B80").Value = ""
$ws.Range("C80").Value = "Remover campo Idioma X ao clicar em [+]"

$ws.Range("C81").Value = "Incluir [-] em Novo Idioma/Instiruição para remover campo"

$ws.Range("C82").Value = 'Corrigir largura do campo "Escreve"'

$ws.Range("B83").Value = "Formação"

$ws.Range("C84").Value = 'Exibir "Órgão financiador" se "Bolsa?" = sim'

$ws.Range("C85").Value = 'Incluir botão [+] em "Órgão Financiador"'

$ws.Range("C86").Value = "Inclruir botão [-] em Novo Órgão para remover campo"

$ws.Range("B87").Value = "UF Nascimento"
$ws.Range("C87").Value = "Aparece se Nacionalidade = Brasil"

# Update view state to match: visible window scrolled up a bit and the
# new block (B79:C87) selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 56
$win.ScrollColumn = 1
[void]$ws.Range("B79:C87").Select()
